$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new row 61 is appended, re-adding the record exactly as it existed
# before the correction below: phone stored as text "09876543" (leading
# zero preserved), no birthday on file, 0 points.
# Clone A60 (still holding the original text value at this point) down to
# A61 so the text type/content is preserved exactly, without touching
# number formats / introducing new cell styles.
$ws.Range("A60").Copy()
$ws.Range("A61").PasteSpecial(-4163)   # xlPasteValues

$ws.Range("B61").Value = ""
$ws.Range("C61").Value = 0

# Row 60: correct the phone number from text "09876543" (leading zero) to
# the numeric value 9876543.
$ws.Range("A60").Value = 9876543
